# "Checklists nominal.xlsx" — Sheet1
# 1) Row 2's Start time / Completion time get moved from 21-Oct-2019 to
#    1-Jan-2019, keeping the same time-of-day.
# 2) The saved selection moves from D13 to C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = (Get-Date -Year 2019 -Month 1 -Day 1 -Hour 15 -Minute 31 -Second 44)
$ws.Range("C2").Value = (Get-Date -Year 2019 -Month 1 -Day 1 -Hour 15 -Minute 32 -Second 58)

[void]$ws.Range("C15").Select()
